# Populate the respondent_first / respondent_last columns (G/H) for the
# rows that were missing them (rows 25-36 of the "discrepancies" sheet).
# These values are simply the first/last name split out of column A's
# "LAST, FIRST[...]" respondent text (same pattern already used for every
# other row on the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ordered list of (row, respondent_first, respondent_last)
$values = @(
    , @(25, "AURORA", "VALENTINE")
    , @(26, "AURORA", "VALENTINE")
    , @(27, "ANGEL", "RODRIGUEZ")
    , @(28, "SAMUEL WILMER", "HARRIS")
    , @(29, "JAMES", "SMITH")
    , @(30, "EDWARD", "EPURE")
    , @(31, "JAMES", "SMITH")
    , @(32, "ZEOLA", "LANCASTER")
    , @(33, "JEFFIE", "BROWN")
    , @(34, "RUDOLPH", "POLSELLI")
    , @(35, "SHERFON", "MCNAIR")
    , @(36, "MAMIE", "COOK")
)

foreach ($entry in $values) {
    $row = $entry[0]
    $first = $entry[1]
    $last = $entry[2]

    # Copy the existing formatting from column F of the same row (already
    # styled like the rest of the table) onto the new G/H cells before
    # writing their values, so the new cells match the sheet's existing
    # look instead of picking up a default style.
    $ws.Range("F$row").Copy()
    $ws.Range("G$row").PasteSpecial(-4122)
    $ws.Range("H$row").PasteSpecial(-4122)

    $ws.Range("G$row").Value = $first
    $ws.Range("H$row").Value = $last
}
